# feat: add 2022-Q3 data
#
# - Insert a new worksheet "2022-Q3" right after "总计" (i.e. before the
#   existing "2022-Q2" sheet) and populate it with the Q3 holdings detail.
# - Update the "总计" summary sheet: the row that used to summarize Q2 now
#   summarizes Q3 (new name + new value), and a new row is appended below
#   it holding the Q2 summary (carrying the value that used to live in the
#   Q2 summary row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) 总计 (summary) sheet updates
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Row 2 becomes the 2022-Q3 summary row.
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("D2").Value = 0.14

# Row 3 is a brand-new row holding what used to be the 2022-Q2 summary.
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.16

# Give A3 the same style as A2 (bold/centered/bordered) without touching
# the value we just wrote.
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Insert the new 2022-Q3 detail sheet right after 总计 (this pushes the
#    existing 2022-Q2 sheet one slot to the right, matching the diff).
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add()
$q3.Name = "2022-Q3"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3.Range("A2").Value = 0

# Columns B:G on the data rows are textual in the source data (codes with
# leading zeros, decimal strings) -- force text so Excel doesn't coerce
# them into numbers (which would e.g. drop the leading zero of "002379").
$q3.Range("B2:G3").NumberFormat = "@"
$q3.Range("B2").Value = "002379"
$q3.Range("C2").Value = "工银瑞信香港中小盘股票（QDII）人民币"
$q3.Range("D2").Value = "1.58"
$q3.Range("E2").Value = "78.58"
$q3.Range("F2").Value = "4.34"
$q3.Range("G2").Value = "0.0686"
$q3.Range("H2").Value = 5

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "002380"
$q3.Range("C3").Value = "工银瑞信香港中小盘股票（QDII）美元"
$q3.Range("D3").Value = "1.58"
$q3.Range("E3").Value = "78.58"
$q3.Range("F3").Value = "4.34"
$q3.Range("G3").Value = "0.0686"
$q3.Range("H3").Value = 5

# Drop the scratch "@" number format again -- the source cells carry no
# explicit style, just the forced-text type.
$q3.Range("B2:G3").ClearFormats()

# Header row + column A use the bold/centered/bordered style (same one
# used on 总计's header/A-column), matching the target sheet.
$summary.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$q3.Range("A2").PasteSpecial(-4122)
$q3.Range("A3").PasteSpecial(-4122)
